$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the casing of the archaeology handbook's filename so it picks up
# the "HARRIS" / uppercase naming convention used for the default manifest row.
$ws.Range("A3").Value = "HANDBOOK_of_archaeology.jpg"

# Widen column A so the (now longer-looking) filenames are fully visible.
# (ColumnWidth is in "characters"; Excel re-quantizes it to pixels on its
# internal grid, so this lands the stored <col width> on ~34.5 - matching
# the manually-resized column A from the original edit.)
$ws.Columns.Item(1).ColumnWidth = 33.7

# Leave the selection where Excel would land after the edit.
$ws.Range("A4").Select()
